$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gpValues = @(10251,2020,10,1,1,59,55,4,1,5,175,131,306,1.34,55,55,39,26,19,8,6,7,4,19,42.1,38.25,16.11,36,51,31,47,2.47,5.88,31.9,17,188.2,87.7,25.74,104.3,7,5,5,5,114,184,204,66.7,51,4,7,36,31,59,1,6,75,183,124,307,1.48,69,44,22,19,26,8,5,5,2,15,53.3,38.38,20.47,36,64,39,39,2.6,4.88,33.3,20.5,186.1,86.3,25.16,79.7,8,8,3,3,120,173,200,65.09999999999999,64,4,8,36,39,44,5,5,62.5)
$gqValues = @(10260,2020,12,1,0,25,66,-41,0,16,160,98,258,1.63,43,59,30,16,16,3,0,6,1,10,30,86,25.8,34,48,22,52,5.2,17.33,17.3,5.8,187.4,86.59999999999999,25.49,97.7,8,5,5,4,128,126,148,57.4,48,6,4,34,22,59,0,0,0,203,100,303,2.03,84,57,38,16,16,10,6,5,1,16,62.5,30.3,18.94,32,48,47,33,2.06,3.3,45.5,30.3,187.2,83.5,24.8,68.59999999999999,11,5,4,2,133,161,195,64.40000000000001,48,17,7,32,47,46,7,6,60)

for ($i = 0; $i -lt $gpValues.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 198).Value = $gpValues[$i]
    $ws.Cells.Item($row, 199).Value = $gqValues[$i]
}
